$d = $word.ActiveDocument
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid wp14"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Changes to Assignment 2 Implementation</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Added </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CropCapabilit</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>y</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">enum </w:t></w:r><w:r><w:t>to</w:t></w:r><w:r><w:t xml:space="preserve"> be used in</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Crop</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Farmer</w:t></w:r><w:r><w:t xml:space="preserve"> to identify if </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Crop</w:t></w:r><w:r><w:t xml:space="preserve"> is ripe or unripe. </w:t></w:r><w:r><w:t>So,</w:t></w:r><w:r><w:t xml:space="preserve"> there is no need to make assumptions if</w:t></w:r><w:r><w:t xml:space="preserve"> the</w:t></w:r><w:r><w:t xml:space="preserve"> crop is ripe or unripe based on the size of </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>allowableAction</w:t></w:r><w:r><w:t xml:space="preserve"> on </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Crop</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Instead of </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CraftWeaponAction</w:t></w:r><w:r><w:t xml:space="preserve"> being responsible </w:t></w:r><w:r><w:t xml:space="preserve">of creating the upgraded weapon from </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ZombieLimb</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t xml:space="preserve">the subclasses of </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ZombieLimb</w:t></w:r><w:r><w:t xml:space="preserve"> are responsib</w:t></w:r><w:r><w:t>le.</w:t></w:r><w:r><w:t xml:space="preserve"> This is done to follow the Open/C</w:t></w:r><w:r><w:t xml:space="preserve">losed Principle. </w:t></w:r><w:r><w:t xml:space="preserve">Now </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CraftWeaponAction</w:t></w:r><w:r><w:t xml:space="preserve"> is open</w:t></w:r><w:r><w:t xml:space="preserve"> for extension by introducing new upgradable items without modifying </w:t></w:r><w:r><w:t xml:space="preserve">it. To achieve this, an </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>upgrade()</w:t></w:r><w:r><w:t xml:space="preserve"> method is added to the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ItemInterface</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">which returns an upgraded </w:t></w:r><w:r><w:t>form of the item if it has one, else, it’ll return null.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Originally, the name of the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ZombieLimb</w:t></w:r><w:r><w:t xml:space="preserve"> is used to identify whether the limb is an arm or leg</w:t></w:r><w:r><w:t>. This is a Connascence of Name</w:t></w:r><w:r><w:t xml:space="preserve"> (CoN)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">that </w:t></w:r><w:r><w:t>an</w:t></w:r><w:r><w:t xml:space="preserve"> IDE cannot pick up and can </w:t></w:r><w:r><w:t xml:space="preserve">lead </w:t></w:r><w:r><w:t>t</w:t></w:r><w:r><w:t xml:space="preserve">o a bug that is hard to identify if the name of the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ZombieLimb</w:t></w:r><w:r><w:t xml:space="preserve"> is changed</w:t></w:r><w:r><w:t xml:space="preserve">. To avoid this, </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ZombieLimb</w:t></w:r><w:r><w:t xml:space="preserve"> is changed to an abstract class and </w:t></w:r><w:r><w:t xml:space="preserve">has two subclasses – </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ZombieArm</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ZombieLeg</w:t></w:r><w:r><w:t xml:space="preserve">. This still has a </w:t></w:r><w:r><w:t>CoN</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>but it works to our benefit as an IDE can easily identify the bug.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Removed </w:t></w:r><w:r><w:t xml:space="preserve">use of </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>instanceof</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>where it</w:t></w:r><w:r><w:t xml:space="preserve"> is not </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">limited to that class </w:t></w:r><w:r><w:t xml:space="preserve">as it is a code smell which restricts polymorphism. It is replaced </w:t></w:r><w:r><w:t xml:space="preserve">with the use of </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Capabilit</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ies</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Added </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>EatCapability</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>enum to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">identify whether an item </w:t></w:r><w:r><w:t>can</w:t></w:r><w:r><w:t xml:space="preserve"> be eaten</w:t></w:r><w:r><w:t xml:space="preserve">, rather than checking if the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Item</w:t></w:r><w:r><w:t xml:space="preserve"> is an instance of </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Food</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Added </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>GroundCapability</w:t></w:r><w:r><w:t xml:space="preserve"> enum to identify whether a ground is sowable</w:t></w:r><w:r><w:t xml:space="preserve">, rather than checking if the </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Ground</w:t></w:r><w:r><w:t xml:space="preserve"> is an instance of </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Dirt</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="11906" w:h="16838"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="708" w:footer="708" w:gutter="0"/><w:cols w:space="708"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng = $d.Range(0, $d.Content.End)
$rng.InsertXML($xml)
